# Commit: "Added affix types for filters"
#
# Adds a new "affix_type" column (column AT / 46) to the Affixes sheet:
#   - AT1 gets the new header label "affix_type" (added to the shared
#     string table)
#   - AT2:AT56 get the numeric flag value 1 for every existing affix row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Affixes")

$ws.Range("AT1").Value = "affix_type"
$ws.Range("AT2:AT56").Value = 1
